$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles only) from the last existing row (68) down to the
# new row (69) so the appended row matches the sheet's established look
# (bold/bordered/centered index column, date-formatted match-date column).
$ws.Range("A68:V68").Copy()
$ws.Range("A69:V69").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's values.
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "south-africa"
$ws.Range("C69").Value = "premier-league"
$ws.Range("D69").Value = "2023-2024"
$ws.Range("E69").Value = 45225.8125
$ws.Range("F69").Value = "Supersport Utd"
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = "Royal AM"
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 1.65
$ws.Range("K69").Value = "26/10/2023 08:22"
$ws.Range("L69").Value = 1.64
$ws.Range("M69").Value = "26/10/2023 11:43"
$ws.Range("N69").Value = 3.35
$ws.Range("O69").Value = "26/10/2023 08:22"
$ws.Range("P69").Value = 3.51
$ws.Range("Q69").Value = "26/10/2023 19:04"
$ws.Range("R69").Value = 5.69
$ws.Range("S69").Value = "26/10/2023 08:22"
$ws.Range("T69").Value = 6.46
$ws.Range("U69").Value = "26/10/2023 19:21"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/south-africa/premier-league/supersport-utd-royal-am/MRrzNLrS/"
